$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Proof-read / amended event-summary outcome texts (column B, rows 15-18).
$ws.Range("B15").Value = "You poured large amounts of  attention and resource into investigating this issue, making sure it was not of ill-intent. One of the neighbouring kingdom caught wind of this and is taking the opportunity to attack you while your guard was down."
$ws.Range("B16").Value = "The results from the poll returned and it further cemented the truth that your villagers did not like you. The neighbouring kingdom took chance of the distress within your kingdom to commit to their attacks."
$ws.Range("B17").Value = "Increasing the taxes had brought about no positive outcomes. In fact, it had caused even more unrest within the kingdom and the villagers started rioting. The neighbouring kingdom took opportunity of this unrest to attack your throne."
$ws.Range("B18").Value = "You managed to quell some of the unrest within the kingdom and raised your ratings. Aware that you were down on resources to defend your kingdom, the neighbouring kingdom decided to launch an attack."
